$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.591.76"
$ws.Range("E2").Value = "  +4.72%  "
$ws.Range("D3").Value = "2.296.25"
$ws.Range("E3").Value = "  +2.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +11.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.568"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.21%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("E9").Value = "  +4.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.40%  "
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.61%  "
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "2.645.16"
$ws.Range("E14").Value = "  +2.87%  "
$ws.Range("D15").Value = "2.295.79"
$ws.Range("E15").Value = "  +2.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.95%  "
$ws.Range("E17").Value = "  +3.80%  "
$ws.Range("D18").Value = "46.559.59"
$ws.Range("E18").Value = "  +5.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.76%  "
$ws.Range("D20").Value = "0.0₃0938"
$ws.Range("E20").Value = "  +3.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.48%  "
$ws.Range("E24").Value = "  +2.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  +3.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.47%  "
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.22%  "
$ws.Range("E31").Value = "  +12.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "146.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.116"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.65%  "
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +19.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.86%  "
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("E44").Value = "  +7.88%  "
$ws.Range("D45").Value = "1.818.14"
$ws.Range("E45").Value = "  +1.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +19.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.195"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "73.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "95.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").Value = "2.520.83"
$ws.Range("E51").Value = "  +2.80%  "
